$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 468.83334
$ws.Range("I53").Value = 402.2
$ws.Range("J53").Value = 579.8889
$ws.Range("K53").Value = 402.2
$ws.Range("L53").Value = 579.8889
$ws.Range("M53").Value = 234.8
$ws.Range("N53").Value = -1853.8889
$ws.Range("H70").Value = 3154.6667
$ws.Range("I70").Value = 2378.6
$ws.Range("J70").Value = 4124.75
$ws.Range("K70").Value = 7135.799999999999
$ws.Range("L70").Value = 12374.25
$ws.Range("M70").Value = -6865.799999999999
$ws.Range("N70").Value = -12914.25
$ws.Range("H73").Value = 3154.6667
$ws.Range("I73").Value = 2378.6
$ws.Range("J73").Value = 4124.75
$ws.Range("K73").Value = 7135.799999999999
$ws.Range("L73").Value = 12374.25
$ws.Range("M73").Value = -6199.799999999999
$ws.Range("N73").Value = -14246.25
$ws.Range("H82").Value = 4597.8
$ws.Range("I82").Value = 4597.8
$ws.Range("K82").Value = 13793.4
$ws.Range("M82").Value = -13387.4
$ws.Range("H85").Value = 4597.8
$ws.Range("I85").Value = 4597.8
$ws.Range("K85").Value = 13793.4
$ws.Range("M85").Value = -12389.4
$ws.Range("H132").Value = 236839.16
$ws.Range("I132").Value = 314804.56
$ws.Range("K132").Value = 944413.6799999999
$ws.Range("M132").Value = -941883.6799999999
$ws.Range("H133").Value = 65875.8
$ws.Range("J133").Value = 65875.8
$ws.Range("L133").Value = 65875.8
$ws.Range("N133").Value = -75995.8
$ws.Range("H135").Value = 3475.8867
$ws.Range("I135").Value = 1393.5122
$ws.Range("J135").Value = 10590.667
$ws.Range("K135").Value = 12541.6098
$ws.Range("L135").Value = 95316.003
$ws.Range("M135").Value = -10006.6098
$ws.Range("N135").Value = -100386.003
$ws.Range("H138").Value = 3544.5881
$ws.Range("I138").Value = 1302.619
$ws.Range("K138").Value = 3907.857
$ws.Range("M138").Value = 1232.143

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 57392.75
$ws.Range("J76").Value = 57392.75
$ws.Range("L76").Value = 57392.75
$ws.Range("N76").Value = -58068.75
$ws.Range("H79").Value = 57392.75
$ws.Range("J79").Value = 57392.75
$ws.Range("L79").Value = 57392.75
$ws.Range("N79").Value = -59732.75
$ws.Range("H132").Value = 1123064.5
$ws.Range("I132").Value = 1379957.1
$ws.Range("J132").Value = 181125
$ws.Range("K132").Value = 4139871.3
$ws.Range("L132").Value = 543375
$ws.Range("M132").Value = -4137341.3
$ws.Range("N132").Value = -548435

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 21046.5
$ws.Range("J88").Value = 21046.5
$ws.Range("L88").Value = 21046.5
$ws.Range("N88").Value = -21858.5
$ws.Range("H91").Value = 21046.5
$ws.Range("J91").Value = 21046.5
$ws.Range("L91").Value = 21046.5
$ws.Range("N91").Value = -23854.5
$ws.Range("H94").Value = 4324.8335
$ws.Range("J94").Value = 5447.625
$ws.Range("L94").Value = 5447.625
$ws.Range("N94").Value = -6349.625
$ws.Range("H134").Value = 1478125.8
$ws.Range("I134").Value = 1792492.4
$ws.Range("K134").Value = 5377477.199999999
$ws.Range("M134").Value = -5374942.199999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 197500
$ws.Range("J88").Value = 197500
$ws.Range("L88").Value = 197500
$ws.Range("N88").Value = -198312
$ws.Range("H91").Value = 197500
$ws.Range("J91").Value = 197500
$ws.Range("L91").Value = 197500
$ws.Range("N91").Value = -200308
$ws.Range("H99").Value = 2927242.8
$ws.Range("I99").Value = 6946810.5
$ws.Range("J99").Value = 3920.818
$ws.Range("K99").Value = 6946810.5
$ws.Range("L99").Value = 3920.818
$ws.Range("M99").Value = -6945312.5
$ws.Range("N99").Value = -6916.818
$ws.Range("H126").Value = 2927242.8
$ws.Range("I126").Value = 6946810.5
$ws.Range("J126").Value = 3920.818
$ws.Range("K126").Value = 20840431.5
$ws.Range("L126").Value = 11762.454
$ws.Range("M126").Value = -20837961.5
$ws.Range("N126").Value = -16702.454
$ws.Range("H132").Value = 6099.72
$ws.Range("I132").Value = 5297.2856
$ws.Range("J132").Value = 10312.5
$ws.Range("K132").Value = 15891.8568
$ws.Range("L132").Value = 30937.5
$ws.Range("M132").Value = -13361.8568
$ws.Range("N132").Value = -35997.5
$ws.Range("H134").Value = 26320002
$ws.Range("I134").Value = 33336678
$ws.Range("J134").Value = 7466.375
$ws.Range("K134").Value = 100010034
$ws.Range("L134").Value = 22399.125
$ws.Range("M134").Value = -100007499
$ws.Range("N134").Value = -27469.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 313423.2
$ws.Range("I113").Value = 555
$ws.Range("J113").Value = 334281.06
$ws.Range("K113").Value = 1665
$ws.Range("L113").Value = 1002843.18
$ws.Range("M113").Value = 505
$ws.Range("N113").Value = -1007183.18
$ws.Range("H131").Value = 34670716
$ws.Range("J131").Value = 25644208
$ws.Range("L131").Value = 76932624
$ws.Range("N131").Value = -76942704
$ws.Range("H137").Value = 1976.8
$ws.Range("J137").Value = 1967.375
$ws.Range("L137").Value = 5902.125
$ws.Range("N137").Value = -16102.125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 21279816
$ws.Range("I132").Value = 26318622
$ws.Range("K132").Value = 78955866
$ws.Range("M132").Value = -78953336
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 31999.5
$ws.Range("J54").Value = 31999.5
$ws.Range("L54").Value = 31999.5
$ws.Range("N54").Value = -33287.5
$ws.Range("H68").Value = 1887.9474
$ws.Range("I68").Value = 1362.9286
$ws.Range("J68").Value = 3358
$ws.Range("K68").Value = 1362.9286
$ws.Range("L68").Value = 3358
$ws.Range("M68").Value = -613.9286
$ws.Range("N68").Value = -4856
$ws.Range("H71").Value = 1887.9474
$ws.Range("I71").Value = 1362.9286
$ws.Range("J71").Value = 3358
$ws.Range("K71").Value = 6814.643
$ws.Range("L71").Value = 16790
$ws.Range("M71").Value = -3070.643
$ws.Range("N71").Value = -24278
$ws.Range("H132").Value = 5058.457
$ws.Range("I132").Value = 3620.5715
$ws.Range("K132").Value = 10861.7145
$ws.Range("M132").Value = -8331.7145

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 33333
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H62").Value = 20499.889
$ws.Range("I62").Value = 19625
$ws.Range("J62").Value = 20749.857
$ws.Range("K62").Value = 19625
$ws.Range("L62").Value = 20749.857
$ws.Range("M62").Value = -19001
$ws.Range("N62").Value = -21997.857
$ws.Range("H63").Value = 39974.5
$ws.Range("J63").Value = 39974.5
$ws.Range("L63").Value = 39974.5
$ws.Range("N63").Value = -41222.5
$ws.Range("H65").Value = 20499.889
$ws.Range("I65").Value = 19625
$ws.Range("J65").Value = 20749.857
$ws.Range("K65").Value = 98125
$ws.Range("L65").Value = 103749.285
$ws.Range("M65").Value = -95005
$ws.Range("N65").Value = -109989.285
$ws.Range("H66").Value = 39974.5
$ws.Range("J66").Value = 39974.5
$ws.Range("L66").Value = 119923.5
$ws.Range("N66").Value = -126163.5
$ws.Range("H80").Value = 300
$ws.Range("J80").Value = 300
$ws.Range("L80").Value = 300
$ws.Range("N80").Value = -2296
$ws.Range("H82").Value = 68979.60000000001
$ws.Range("J82").Value = 68979.60000000001
$ws.Range("L82").Value = 68979.60000000001
$ws.Range("N82").Value = -69745.60000000001
$ws.Range("H83").Value = 300
$ws.Range("J83").Value = 300
$ws.Range("L83").Value = 900
$ws.Range("N83").Value = -10884
$ws.Range("H85").Value = 68979.60000000001
$ws.Range("J85").Value = 68979.60000000001
$ws.Range("L85").Value = 68979.60000000001
$ws.Range("N85").Value = -71631.60000000001
$ws.Range("H126").Value = 7430.5
$ws.Range("I126").Value = 7055.4
$ws.Range("K126").Value = 21166.2
$ws.Range("M126").Value = -18696.2
$ws.Range("H132").Value = 5934.4194
$ws.Range("I132").Value = 5312.25
$ws.Range("J132").Value = 7065.636
$ws.Range("K132").Value = 15936.75
$ws.Range("L132").Value = 21196.908
$ws.Range("M132").Value = -13406.75
$ws.Range("N132").Value = -26256.908
